# Update the second onboarding row ("Sheet2") with a fresh set of
# auto-generated user details (firstName / lastName / emailID / mobNumber)
# as part of the 100-user onboarding flow test-data refresh.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A2").Value = "Gary"
$ws2.Range("B2").Value = "Torphy"
$ws2.Range("C2").Value = "66S6O@mailinator.com"
$ws2.Range("D2").Value = "9826031438"
